# The workbook used to ship with two sheets:
#   - "Sheet1": an empty leftover sheet
#   - "Sheet2": the actual BoM data (with its query table / external data link)
# Clean this up: drop the empty sheet and rename the data sheet to "Sheet1".

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()
$excel.DisplayAlerts = $true

# Re-fetch by name after the delete (don't reuse a reference captured
# before the delete, since sheet indices shift).
$wb.Worksheets.Item("Sheet2").Name = "Sheet1"
$wb.Worksheets.Item("Sheet1").Activate()
